$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Temp")
$ws.Activate()

$ws.Range("A3").Value = "tem002"
$ws.Range("B3").Value = "temp detail rate hotel API"

$ws.Range("A4").Select()
